$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 123 (shifts the existing rows 123..146 down to 124..147)
$ws.Rows.Item(123).Insert()

# Populate the newly inserted row 123 with the new weekly record
$ws.Range("A123").Value = 5
$ws.Range("B123").Value = "Macroferia Regional de Talca"
$ws.Range("C123").Value = "Maule"
$ws.Range("D123").Value = 44504
$ws.Range("E123").Value = 7
$ws.Range("F123").Value = 100112045
$ws.Range("G123").Value = "Zapallo"
$ws.Range("H123").Value = "Paine"
$ws.Range("I123").Value = "1a (guarda)"
$ws.Range("J123").Value = 2500
$ws.Range("K123").Value = 80
$ws.Range("L123").Value = 80
$ws.Range("M123").Value = 80
$ws.Range("N123").Value = '$/kilo (volumen en unidades)'
$ws.Range("O123").Value = "Región del Maule"
$ws.Range("P123").Value = 80
$ws.Range("Q123").Value = 1
$ws.Range("R123").Value = "Hortaliza"
